$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New password values for column D (rows 2-31)
$passwords = @{
    2  = "ZnQyx4P"
    3  = "AtfsIpl"
    4  = "mpSzJv3"
    5  = "jCYAttN"
    6  = "tgBpxYb"
    7  = "RrICFED"
    8  = "JQlQ9NZ"
    9  = "DqR6KB2"
    10 = "A941tdW"
    11 = "ejZr2dX"
    12 = "n4D3AXE"
    13 = "JZ3mNbv"
    14 = "TPFGN1T"
    15 = "rI1rheM"
    16 = "GemE5BE"
    17 = "Cfj6HOb"
    18 = "d6qU3j1"
    19 = "n0R7veo"
    20 = "MRfvbL0"
    21 = "XnDrjGb"
    22 = "OlJ8vrp"
    23 = "KqpF2OI"
    24 = "AJB43FK"
    25 = "Xm63qfG"
    26 = "CkoJR9C"
    27 = "CYflJNi"
    28 = "kYDMPTY"
    29 = "GajbI0k"
    30 = "VYvG9Ai"
    31 = "1S59RT0"
}

foreach ($row in $passwords.Keys) {
    $ws.Cells.Item($row, 4).Value = $passwords[$row]
}

# Column B (Class) rows 8-31 become lowercase "a" (rows 2-7 remain "A")
for ($row = 8; $row -le 31; $row++) {
    $ws.Cells.Item($row, 2).Value = "a"
}
